# Combustuble_A.xlsx edit script
# - Renames "Gasvión" -> "Gasavión" everywhere (table headers, formulas).
# - Converts the numeric "Mes" column (1-12) into Spanish month abbreviations
#   ("Ene.", "Feb.", ... "Dic.") stored as text, in the same order the
#   abbreviations are first encountered top-to-bottom (Ago. first).
# - Fixes two typos in footnote text (title + source note).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Convert the "Mes" column (column C, rows 6-85) from numbers to text
#    month abbreviations. Do this FIRST so the new shared strings for the
#    month abbreviations are created (in Ago..Sep order) before the other
#    renamed strings below.
# ---------------------------------------------------------------------
$monthAbbrev = @{
    1  = "Ene."
    2  = "Feb."
    3  = "Mar."
    4  = "Abr."
    5  = "May."
    6  = "Jun."
    7  = "Jul."
    8  = "Ago."
    9  = "Sep."
    10 = "Oct."
    11 = "Nov."
    12 = "Dic."
}

for ($r = 6; $r -le 85; $r++) {
    $monthNum = $ws.Cells.Item($r, 3).Value2
    $abbrev = $monthAbbrev[[int]$monthNum]
    $ws.Cells.Item($r, 3).Value2 = $abbrev
}

# ---------------------------------------------------------------------
# 2) Fix the report title typo ("aeropuerts" -> "aeropuertos").
# ---------------------------------------------------------------------
$ws.Range("B2").Value2 = "Consumo de combustible en los aeropuertos de ASA y grupo aeroportuarios"

# ---------------------------------------------------------------------
# 3) Rename "Gasvión" -> "Gasavión" in the table headers (this renames the
#    underlying Excel Table column names as well).
# ---------------------------------------------------------------------
$ws.Range("F5").Value2 = "Gasavión"
$ws.Range("K5").Value2 = "Gasavión  "

# ---------------------------------------------------------------------
# 4) Update the calculated formulas that reference the "Gasvión" / "Gasvión  "
#    table columns so they reference the renamed "Gasavión" / "Gasavión  "
#    columns instead (column D = Total, column I = Total  ).
# ---------------------------------------------------------------------
for ($r = 6; $r -le 85; $r++) {
    $ws.Cells.Item($r, 4).Formula = "=SUM(Tabla3[[#This Row],[Turbosina]:[Gasavión]])"
}

for ($r = 6; $r -le 32; $r++) {
    $ws.Cells.Item($r, 9).Formula = "=SUM(Tabla3[[#This Row],[Turbosina  ]:[Gasavión  ]])"
}

# ---------------------------------------------------------------------
# 5) Fix the source footnote (add trailing period after the URL).
# ---------------------------------------------------------------------
$ws.Range("B94").Value2 = "Fuente: ASA. Aeropuertos y Servicios Auxiliares, litros y servicios por venta de combustibles. En: www.asa.gob.mx."
